$wb = $excel.ActiveWorkbook

# --- Overview sheet: bump the "Latest HO Xliff Generate Date" shared by
#     rows 4-7 (5dfa90da, 75eeee6c, 8bced00c, dca63bb7) from 12:26:23 ->
#     12:26:37. These rows (and the de-de sheet's "Latest Handoff
#     Datetime" column for the same rows) all share one string value, so
#     all of them need to be re-written to keep sharing the new text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-15 12:26:37"
$wsOverview.Range("G5").Value = "2016-08-15 12:26:37"
$wsOverview.Range("G6").Value = "2016-08-15 12:26:37"
$wsOverview.Range("G7").Value = "2016-08-15 12:26:37"

# --- zh-cn sheet: rows 4-7 (the "Ready for handoff" rows) move from
#     Priority "low" -> "ht", and their handoff datetime advances from
#     12:26:17 -> 12:26:32.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("E7").Value = "ht"

$wsZhCn.Range("H4").Value = "2016-08-15 12:26:32"
$wsZhCn.Range("H5").Value = "2016-08-15 12:26:32"
$wsZhCn.Range("H6").Value = "2016-08-15 12:26:32"
$wsZhCn.Range("H7").Value = "2016-08-15 12:26:32"

# --- de-de sheet: rows 4-7 Priority also moves "low" -> "ht", and their
#     "Latest Handoff Datetime" shares the same string as Overview!G4:G7
#     above, so it must track the same new value.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("E7").Value = "ht"

$wsDeDe.Range("H4").Value = "2016-08-15 12:26:37"
$wsDeDe.Range("H5").Value = "2016-08-15 12:26:37"
$wsDeDe.Range("H6").Value = "2016-08-15 12:26:37"
$wsDeDe.Range("H7").Value = "2016-08-15 12:26:37"
